$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 (PPF004 / Swathipriya): the Week01 task is now completed, so the
# row is restyled to match the other "Completed" rows (e.g. row 25): the
# student-name and completion-status cells get the green "Completed" fill,
# the pending-task note is cleared, status flips to Completed, and the
# streak increments.
$ws.Range("B25").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Range("E24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("E24").Value = "Completed"
$ws.Range("D24").Value = ""
$ws.Range("F24").Value = 1

$ws.Range("D24").Select() | Out-Null
